$d = $word.ActiveDocument

function Bold-SubRange($baseStart, $fullText, $subText) {
    $idx = $fullText.IndexOf($subText)
    if ($idx -ge 0) {
        $s = $baseStart + $idx
        $e = $s + $subText.Length
        $r = $d.Range($s, $e)
        $r.Bold = $true
        $r.BoldBi = $true
    }
}

function Replace-WithinFound($findText, $newText, $boldParts) {
    $rng = $d.Content
    $find = $rng.Find
    $find.ClearFormatting()
    $found = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $start = $rng.Start
        $rng.Text = $newText
        foreach ($b in $boldParts) {
            Bold-SubRange $start $newText $b
        }
    }
    return $found
}

# ----------------------------------------------------------------------
# 1) "O sistema exibe a listagem ..." -> "O sistema apresenta a listagem ..."
#    (bold "apresenta a listagem")
# ----------------------------------------------------------------------
Replace-WithinFound `
    "O sistema exibe a listagem de todas as reservas, permitindo filtrar por área, data e " `
    "O sistema apresenta a listagem de todas as reservas, permitindo filtrar por área, data e " `
    @("apresenta a listagem")

# ----------------------------------------------------------------------
# 2) "Faz a edição com espaços disponíveis." ->
#    "O ator atualiza os dados da reserva (data/horário) conforme disponibilidade.."
#    (bold "atualiza os dados")
# ----------------------------------------------------------------------
Replace-WithinFound `
    "Faz a edição com espaços disponíveis." `
    "O ator atualiza os dados da reserva (data/horário) conforme disponibilidade.." `
    @("atualiza os dados")

# ----------------------------------------------------------------------
# 3) "Tela_004" (first, "reserva desejada (...)").-> "Ver_Wireframe_004"
# ----------------------------------------------------------------------
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$found = $find.Execute("desejada (Tela_004)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $subRng = $d.Range($rng.Start, $rng.End)
    $subFind = $subRng.Find
    $subFind.ClearFormatting()
    $subFind.Replacement.ClearFormatting()
    $subFind.Execute("Tela_", $true, $false, $false, $false, $false, $true, 1, $false, "Ver_Wireframe_", 1) | Out-Null
}

# ----------------------------------------------------------------------
# 4) "Ator seleciona a opção remover reserva ." -> left as is (text identical);
#    only proofErr markup differs in the source diff, no visible change.
# ----------------------------------------------------------------------

# ----------------------------------------------------------------------
# 5) "Sistema atualiza a tela de reservas " -> "O sistema atualiza a listagem de reservas "
#    (bold "atualiza a listagem")
# ----------------------------------------------------------------------
Replace-WithinFound `
    "Sistema atualiza a tela de reservas " `
    "O sistema atualiza a listagem de reservas " `
    @("atualiza a listagem")

# ----------------------------------------------------------------------
# 6) "Tela_004" (second, "...de reservas (...)") -> "Ver_Wireframe_004"
# ----------------------------------------------------------------------
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$found = $find.Execute("de reservas (Tela_004)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $subRng = $d.Range($rng.Start, $rng.End)
    $subFind = $subRng.Find
    $subFind.ClearFormatting()
    $subFind.Replacement.ClearFormatting()
    $subFind.Execute("Tela_", $true, $false, $false, $false, $false, $true, 1, $false, "Ver_Wireframe_", 1) | Out-Null
}

# ----------------------------------------------------------------------
# 7) ": O ator pressiona botão Cancelar. O sistema retorna ao passo 2 da " ->
#    ": O ator cancela a operação. O sistema apresenta novamente a listagem de reservas. O sistema retorna ao passo 2 da "
#    (bold "cancela a operação" and "apresenta novamente")
# ----------------------------------------------------------------------
Replace-WithinFound `
    ": O ator pressiona botão Cancelar. O sistema retorna ao passo 2 da " `
    ": O ator cancela a operação. O sistema apresenta novamente a listagem de reservas. O sistema retorna ao passo 2 da " `
    @("cancela a operação", "apresenta novamente")

# ----------------------------------------------------------------------
# 8) "Tela_" -> "Ver_Wireframe_" ("Visualizar Reserva na (Tela_004)")
# ----------------------------------------------------------------------
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$found = $find.Execute("Visualizar Reserva na (Tela_004)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $subRng = $d.Range($rng.Start, $rng.End)
    $subFind = $subRng.Find
    $subFind.ClearFormatting()
    $subFind.Replacement.ClearFormatting()
    $subFind.Execute("Tela_", $true, $false, $false, $false, $false, $true, 1, $false, "Ver_Wireframe_", 1) | Out-Null
}

# ----------------------------------------------------------------------
# 9) "Tela_" -> "Ver_Wireframe_" ("selecionada (Tela_006)")
# ----------------------------------------------------------------------
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$found = $find.Execute("selecionada (Tela_006)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $subRng = $d.Range($rng.Start, $rng.End)
    $subFind = $subRng.Find
    $subFind.ClearFormatting()
    $subFind.Replacement.ClearFormatting()
    $subFind.Execute("Tela_", $true, $false, $false, $false, $false, $true, 1, $false, "Ver_Wireframe_", 1) | Out-Null
}
